# Generate Report for Handoff
#
# Refreshes the localization-status report after a new handoff package was
# generated: the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
# move forward a few seconds/minutes, and the newly-handed-off rows get their
# Priority marked "ht" (handoff type) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 12, 13, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G" + $r).Value = "2016-08-27 14:22:30"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H" + $r).Value = "2016-08-27 14:22:26"
    $zhcn.Range("E" + $r).Value = "ht"
}

# --- de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H" + $r).Value = "2016-08-27 14:22:30"
    $dede.Range("E" + $r).Value = "ht"
}
